$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet by duplicating the "2021-Q4" sheet (same
#    fund-holdings layout/styling) and dropping the copy right before "总计".
# ---------------------------------------------------------------------------
$srcQ4 = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")
$srcQ4.Copy($total)

# The duplicate lands immediately before "总计"; Excel auto-names it
# "2021-Q4 (2)" - grab it by that name (re-querying avoids relying on a
# stale Index captured before the sheet got inserted) and rename it.
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"
$total = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------------
# 2. Overwrite the copied data with the 2022-Q1 fund holdings.
#    Numeric-looking values (fund code / percentages) must stay plain text,
#    matching the source data, so they are staged through a Text-formatted
#    helper cell and pasted in as values - this avoids Excel silently
#    re-interpreting e.g. "001643" or "0.4080" as numbers (which would drop
#    the leading zero / trailing zero).
# ---------------------------------------------------------------------------
$stage = $newSheet.Range("Z1")
$stage.NumberFormat = "@"

function Set-TextValue($cell, $text) {
    $stage.Value = $text
    $stage.Copy()
    $cell.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

Set-TextValue $newSheet.Cells.Item(2, 2) "001643"
$newSheet.Cells.Item(2, 3).Value = "汇丰晋信智造先锋股票A"
Set-TextValue $newSheet.Cells.Item(2, 4) "29.09"
Set-TextValue $newSheet.Cells.Item(2, 5) "92.99"
Set-TextValue $newSheet.Cells.Item(2, 6) "3.74"
Set-TextValue $newSheet.Cells.Item(2, 7) "1.0880"
$newSheet.Cells.Item(2, 8).Value = 8

Set-TextValue $newSheet.Cells.Item(3, 2) "001644"
$newSheet.Cells.Item(3, 3).Value = "汇丰晋信智造先锋股票C"
Set-TextValue $newSheet.Cells.Item(3, 4) "10.91"
Set-TextValue $newSheet.Cells.Item(3, 5) "92.99"
Set-TextValue $newSheet.Cells.Item(3, 6) "3.74"
Set-TextValue $newSheet.Cells.Item(3, 7) "0.4080"
$newSheet.Cells.Item(3, 8).Value = 8

$stage.Clear()
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Update the "总计" summary sheet: add a new first data row for 2022-Q1
#    and push the existing rows down, keeping the running index in column A.
# ---------------------------------------------------------------------------
$total.Cells.Item(6, 1).Value = 4
$total.Range("A5").Copy()
$total.Range("A6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats (reuse A-col style)
$total.Cells.Item(6, 2).Value = "2021-Q1"
$total.Cells.Item(6, 3).Value = 6
$total.Cells.Item(6, 4).Value = 1.44

$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(5, 2).Value = "2021-Q2"
$total.Cells.Item(5, 3).Value = 2
$total.Cells.Item(5, 4).Value = 1.59

$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = "2021-Q3"
$total.Cells.Item(4, 3).Value = 2
$total.Cells.Item(4, 4).Value = 1.65

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2021-Q4"
$total.Cells.Item(3, 3).Value = 2
$total.Cells.Item(3, 4).Value = 0.83

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 2
$total.Cells.Item(2, 4).Value = 1.5

$excel.CutCopyMode = $false
